$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("E1").Value = "Act-Avg"

# Row 2 (literal formulas, not shared)
$ws.Range("D2").Formula = "=(C2-B2)^2"
$ws.Range("E2").Formula = "=(B2-B`$8)^2"

# Rows 3:6 (shared formulas)
$ws.Range("D3:D6").Formula = "=(C3-B3)^2"
$ws.Range("E3:E6").Formula = "=(B3-B`$8)^2"

# Row 8: Avg / RSS / TSS summary row
$ws.Range("A8").Value = "Avg"
$ws.Range("B8").Formula = "=AVERAGE(B2:B6)"
$ws.Range("C8").Value = "RSS"
$ws.Range("D8").Formula = "=SUM(D2:D6)"
$ws.Range("F8").Value = "TSS"
$ws.Range("G8").Formula = "=SUM(E2:E6)"

# Row 10: R2
$ws.Range("C10").Value = "R2"
$ws.Range("D10").Formula = "=1-(D8/G8)"

# Update selection to match target (E10 active cell)
$ws.Range("E10").Select() | Out-Null
